$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph - its start is where deletion begins.
$r1 = $d.Content
$r1.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Expand(4) | Out-Null
$startP = $r1.Start

# Locate the "... Creative Commons Attribution" paragraph (the copyright line).
$r2 = $d.Content
$r2.Find.Execute("Creative Commons Attribution", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Expand(4) | Out-Null

# Also sweep up the blank paragraph that immediately follows the copyright line.
$afterRange = $d.Range($r2.End, $r2.End)
$afterRange.Expand(4) | Out-Null
$endP = $afterRange.End

# Delete the "Ver no Jupiter ..." paragraph, the copyright paragraph, and the
# trailing blank paragraph in one shot (leaves the existing blank paragraph
# right after "LOB1004: ..." untouched).
$full = $d.Range($startP, $endP)
$full.Delete() | Out-Null
